$d = $word.ActiveDocument

# Find the paragraph that holds the "in progress" text (style "FirstParagraph")
# sitting right after the "F2025" date paragraph and remove it entirely,
# including its trailing paragraph mark.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "in progress") {
        $p.Range.Delete()
    }
}
